$d = $word.ActiveDocument

$replacements = @(
    @("53×80=", "94×79="),
    @("78×90=", "57×23="),
    @("87×88=", "88×62="),
    @("51×93=", "97×54="),
    @("16×78=", "51×84="),
    @("53×53=", "67×81="),
    @("82×21=", "57×48="),
    @("98×48=", "73×76="),
    @("57×82=", "95×86="),
    @("36×26=", "28×40="),
    @("67×37=", "27×59="),
    @("52×34=", "19×58="),
    @("53×35=", "93×31="),
    @("50×61=", "30×18="),
    @("25×35=", "57×60="),
    @("15×22=", "60×81="),
    @("38×13=", "43×92="),
    @("69×47=", "64×66="),
    @("46×41=", "99×64="),
    @("40×95=", "40×19="),
    @("55×26=", "95×87="),
    @("62×31=", "39×21="),
    @("42×26=", "84×55="),
    @("27×40=", "83×37="),
    @("41×54=", "84×35=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
